# 22/10/2017 MAMATHA CHICK IN
#
# (1) The "... THU Oct 19" / " 11:17:48 PDT 2017" timestamp was stored as
#     two adjacent runs; collapse it back into a single run by finding the
#     full text and "replacing" it with itself.
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    " THU Oct 19 11:17:48 PDT 2017", $true, $false, $false, $false, $false,
    $true, 1, $false, " THU Oct 19 11:17:48 PDT 2017", 2) | Out-Null

# (2) Append a new TSNP purchase record (SAT Oct 21) right after the very
#     last existing "Amount Received mode ... - CASH" line in the document.

# Find the last paragraph whose text is "Amount Received mode" + 2 tabs +
# "- CASH" -- several such blocks exist earlier in the document, we want
# the final one.
$anchorIdx = -1
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $paraText = $d.Paragraphs($i).Range.Text
    if ($paraText -eq "Amount Received mode`t`t- CASH`r") {
        $anchorIdx = $i
        break
    }
}

# Insert a brand-new empty paragraph right after it, then fill each new
# paragraph in turn, always inserting the *next* paragraph right after the
# one we just wrote.
$idx = $anchorIdx
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
# paragraph $idx is a new, blank separator paragraph -- leave it empty.

$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1

# -- SAT Oct 21  11:24:24 PDT 2017  (two runs: date, then time) --
$d.Paragraphs($idx).Range.Text = "SAT Oct 21"
$p = $d.Paragraphs($idx).Range
$insPos = $p.End - 1
$d.Range($insPos, $insPos).InsertAfter(" 11:24:24 PDT 2017")

$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.Text = "Person Name`t`t`t`t- TSNP"

$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.Text = "---------------------------------------------------------------"

$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.Text = "Item Name`t`t`t`t- CARROT EVE"

$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.Text = "Number of Pockets`t`t`t- 1"

$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.Text = "Number of KGs`t`t`t- 91"

$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.Text = "Rate`t`t`t`t`t- 50"

$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.Text = "Total Price`t`t`t`t- 4550.0"

$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.Text = "Amount balance`t`t`t- 17878.0"
$d.Paragraphs($idx).Range.Font.Bold = 1

# Two trailing blank "Plain Text" paragraphs (not bold -- the bold
# paragraph mark above must not leak into them).
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.Font.Bold = 0

$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.Font.Bold = 0

Write-Output "Inserted new TSNP record; document now has $($d.Paragraphs.Count) paragraphs."
